$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.908.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").Value = "'1.642.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.23%  "

$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "'215.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").Value = "'1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.35%  "

$ws.Range("D8").Value = "'0.2572"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "'0.06402"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("D10").Value = "'19.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.20%  "

$ws.Range("D11").Value = "'0.07797"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.07%  "

$ws.Range("D12").Value = "'4.287"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.18%  "

$ws.Range("D13").Value = "'1.650.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.11%  "

$ws.Range("D14").Value = "'0.5442"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.05%  "

$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("D16").Value = "'64.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.98%  "

$ws.Range("D17").Value = "'25.949.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").Value = "'1.006"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.26%  "

$ws.Range("D19").Value = "'198.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.26%  "

$ws.Range("D20").Value = "'4.400"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.62%  "

$ws.Range("D21").Value = "'9.982"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").Value = "'6.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.03%  "

$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").Value = "'1.874"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.70%  "

$ws.Range("D25").Value = "'140.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.52%  "

$ws.Range("D26").Value = "'0.1145"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.59%  "

$ws.Range("D27").Value = "'6.857"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.05%  "

$ws.Range("D28").Value = "'15.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("E29").Value = "  +0.42%  "

$ws.Range("D30").Value = "'0.04963"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.67%  "

$ws.Range("D31").Value = "'3.268"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.59%  "

$ws.Range("D32").Value = "'3.198"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.18%  "

$ws.Range("D33").Value = "'1.535"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.25%  "

$ws.Range("D34").Value = "'2.373"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.47%  "

$ws.Range("D35").Value = "'0.8951"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.79%  "

$ws.Range("D36").Value = "'2.612"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.11%  "

$ws.Range("D37").Value = "'1.141.71"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("D38").Value = "'0.5549"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.06%  "

$ws.Range("D39").Value = "'0.01563"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "

$ws.Range("D40").Value = "'1.006"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.19%  "

$ws.Range("D41").Value = "'5.678"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.46%  "

$ws.Range("D42").Value = "'0.8220"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.71%  "

$ws.Range("D43").Value = "'99.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.23%  "

$ws.Range("D44").Value = "'0.0₈122"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.20%  "

$ws.Range("D45").Value = "'1.777.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("D46").Value = "'0.4527"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").Value = "'55.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.05%  "

$ws.Range("D48").Value = "'1.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.35%  "

$ws.Range("D49").Value = "'0.05059"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.51%  "

$ws.Range("D50").Value = "'1.006"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.12%  "

$ws.Range("D51").Value = "'0.09522"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.42%  "
